$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $r = $ws.Range($rangeAddr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue "D2" '89.633.63'
$ws.Range("E2").Value = '  +0.52%  '
Set-TextValue "D3" '3.044.81'
$ws.Range("E3").Value = '  -3.56%  '
Set-TextValue "D4" '1.00'
$ws.Range("E4").Value = '  -0.03%  '
Set-TextValue "D5" '210.57'
$ws.Range("E5").Value = '  -2.43%  '
Set-TextValue "D6" '612.17'
$ws.Range("E6").Value = '  -3.32%  '
Set-TextValue "D7" '0.362'
$ws.Range("E7").Value = '  -9.11%  '
Set-TextValue "D8" '0.878'
$ws.Range("E8").Value = '  +19.09%  '
$ws.Range("E9").Value = '  +0.00%  '
Set-TextValue "D10" '3.042.80'
$ws.Range("E10").Value = '  -3.55%  '
Set-TextValue "D11" '0.662'
$ws.Range("E11").Value = '  +19.09%  '
$ws.Range("E12").Value = '  +4.71%  '
Set-TextValue "D13" '0.0000238'
$ws.Range("E13").Value = '  -5.73%  '
Set-TextValue "D14" '5.36'
$ws.Range("E14").Value = '  +0.76%  '
Set-TextValue "D15" '88.605.65'
$ws.Range("E15").Value = '  -0.40%  '
Set-TextValue "D16" '32.01'
$ws.Range("E16").Value = '  -1.62%  '
$ws.Range("E17").Value = '  -3.90%  '
Set-TextValue "D18" '3.025.56'
$ws.Range("E18").Value = '  -3.65%  '
Set-TextValue "D19" '3.35'
$ws.Range("E19").Value = '  -1.46%  '
Set-TextValue "D20" '0.0000215'
$ws.Range("E20").Value = '  -7.65%  '
Set-TextValue "D21" '13.36'
$ws.Range("E21").Value = '  +0.15%  '
Set-TextValue "D22" '425.61'
$ws.Range("E22").Value = '  -0.52%  '
Set-TextValue "D23" '5.01'
$ws.Range("E23").Value = '  +2.11%  '
Set-TextValue "D24" '8.20'
$ws.Range("E24").Value = '  -2.62%  '
Set-TextValue "D25" '5.39'
$ws.Range("E25").Value = '  -0.74%  '
Set-TextValue "D26" '83.60'
$ws.Range("E26").Value = '  +2.81%  '
Set-TextValue "D27" '11.65'
$ws.Range("E27").Value = '  +0.72%  '
$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D28" '1.00'
$ws.Range("E28").Value = '  -0.63%  '
$ws.Range("B29").Value = 'Cronos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D29" '0.162'
$ws.Range("E29").Value = '  +1.19%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue "D30" '1.01'
$ws.Range("E30").Value = '  +6.01%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue "D31" '8.18'
$ws.Range("E31").Value = '  -0.80%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue "D32" '503.91'
$ws.Range("E32").Value = '  -1.69%  '
$ws.Range("B33").Value = 'dogwifhat'
$ws.Range("C33").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D33" '3.68'
$ws.Range("E33").Value = '  -9.41%  '
$ws.Range("B34").Value = 'RenderToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue "D34" '6.63'
$ws.Range("E34").Value = '  -6.92%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D35" '22.86'
$ws.Range("E35").Value = '  +3.97%  '
$ws.Range("B36").Value = 'PancakeSwap'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D36" '1.79'
$ws.Range("E36").Value = '  -2.88%  '
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue "D37" '1.23'
$ws.Range("E37").Value = '  -5.23%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D38" '0.131'
$ws.Range("E38").Value = '  -6.49%  '
$ws.Range("B39").Value = 'WhiteBITCoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue "D39" '22.23'
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue "D40" '1.00'
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("B42").Value = 'PolygonEcosystemToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue "D42" '0.362'
$ws.Range("E42").Value = '  -0.99%  '
$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D43" '0.137'
$ws.Range("E43").Value = '  +8.82%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D44" '1.83'
$ws.Range("E44").Value = '  -2.48%  '
$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D45" '145.81'
$ws.Range("E45").Value = '  -0.05%  '
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D46" '43.30'
$ws.Range("E46").Value = '  -1.28%  '
$ws.Range("B47").Value = 'Hedera'
$ws.Range("C47").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D47" '0.0681'
$ws.Range("E47").Value = '  +8.49%  '
$ws.Range("B48").Value = 'Filecoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D48" '4.07'
$ws.Range("E48").Value = '  +3.15%  '
$ws.Range("B49").Value = 'ImmutableX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D49" '1.21'
$ws.Range("E49").Value = '  +1.95%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D50" '159.26'
$ws.Range("E50").Value = '  -4.19%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue "D51" '0.696'
$ws.Range("E51").Value = '  -4.11%  '
